$wb = $excel.ActiveWorkbook

# The stock sheet ("股票") is the 5th sheet in the workbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column at I (pushing old date/legislator_name/legislator_id
# columns from I,J,K to J,K,L) and append two new trailing columns at M,N
# for source_file / index. Column-insert carries the neighbouring column's
# style across automatically, matching the header (s=1) / data (s=2) styles.
$ws.Columns("I:I").Insert()
$ws.Columns("M:N").Insert()

# Header row (row 1)
$ws.Range("I1").Value = "category"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Data row (row 2)
$ws.Range("I2").Value = "normal"
$ws.Range("M2").Value = "tmp38461"
$ws.Range("N2").Value = 59
